# Auto-generated edit script: updates cached market-price / profit values
# in the Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) to reflect
# a refreshed Universalis market-data pull, per the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 232.57143
$ws.Range("I2").Value = 91.8
$ws.Range("J2").Value = 584.5
$ws.Range("K2").Value = 91.8
$ws.Range("L2").Value = 584.5
$ws.Range("M2").Value = 21.2
$ws.Range("N2").Value = -810.5

$ws.Range("H28").Value = 1908.9333
$ws.Range("I28").Value = 1268.6364
$ws.Range("J28").Value = 3669.75
$ws.Range("K28").Value = 1268.6364
$ws.Range("L28").Value = 3669.75
$ws.Range("M28").Value = -783.6364000000001
$ws.Range("N28").Value = -4639.75

$ws.Range("H76").Value = 2928.4285
$ws.Range("I76").Value = 2928.4285
$ws.Range("K76").Value = 2928.4285
$ws.Range("M76").Value = -2613.4285

$ws.Range("H79").Value = 2928.4285
$ws.Range("I79").Value = 2928.4285
$ws.Range("K79").Value = 2928.4285
$ws.Range("M79").Value = -1836.4285

$ws.Range("H107").Value = 713.2857
$ws.Range("I107").Value = 713.2857
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 713.2857
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1206.7143
$ws.Range("N107").ClearContents()

$ws.Range("H135").Value = 25000830
$ws.Range("I135").Value = 788.3514
$ws.Range("J135").Value = 333334700
$ws.Range("K135").Value = 7095.1626
$ws.Range("L135").Value = 3000012300
$ws.Range("M135").Value = -4560.1626
$ws.Range("N135").Value = -3000017370

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6391.3125
$ws.Range("I32").Value = 4378.615
$ws.Range("J32").Value = 15113
$ws.Range("K32").Value = 4378.615
$ws.Range("L32").Value = 15113
$ws.Range("M32").Value = -4091.615
$ws.Range("N32").Value = -15687

$ws.Range("H52").Value = 19000
$ws.Range("J52").Value = 19000
$ws.Range("L52").Value = 19000
$ws.Range("N52").Value = -19636

$ws.Range("H92").Value = 25748.75
$ws.Range("J92").Value = 25748.75
$ws.Range("L92").Value = 25748.75
$ws.Range("N92").Value = -30740.75

$ws.Range("H110").Value = 2135.1428
$ws.Range("I110").Value = 2362.375
$ws.Range("J110").Value = 1832.1666
$ws.Range("K110").Value = 2362.375
$ws.Range("L110").Value = 1832.1666
$ws.Range("M110").Value = -317.375
$ws.Range("N110").Value = -5922.1666

$ws.Range("H137").Value = 38333.168
$ws.Range("J137").Value = 38333.168
$ws.Range("L137").Value = 38333.168
$ws.Range("N137").Value = -48533.168

$ws.Range("H141").Value = 33191.8
$ws.Range("J141").Value = 33191.8
$ws.Range("L141").Value = 33191.8
$ws.Range("N141").Value = -43551.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29115.6
$ws.Range("J55").Value = 29115.6
$ws.Range("L55").Value = 29115.6
$ws.Range("N55").Value = -29661.6

$ws.Range("H86").Value = 1827.1666
$ws.Range("I86").Value = 1638.3529
$ws.Range("J86").Value = 2285.7144
$ws.Range("K86").Value = 1638.3529
$ws.Range("L86").Value = 2285.7144
$ws.Range("M86").Value = -515.3529000000001
$ws.Range("N86").Value = -4531.7144

$ws.Range("H89").Value = 1827.1666
$ws.Range("I89").Value = 1638.3529
$ws.Range("J89").Value = 2285.7144
$ws.Range("K89").Value = 8191.7645
$ws.Range("L89").Value = 11428.572
$ws.Range("M89").Value = -2575.7645
$ws.Range("N89").Value = -22660.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 17000
$ws.Range("J59").Value = 17000
$ws.Range("L59").Value = 17000
$ws.Range("N59").Value = -19290

$ws.Range("H62").Value = 2322.3447
$ws.Range("I62").Value = 2288.5454
$ws.Range("J62").Value = 2428.5715
$ws.Range("K62").Value = 2288.5454
$ws.Range("L62").Value = 2428.5715
$ws.Range("M62").Value = -1664.5454
$ws.Range("N62").Value = -3676.5715

$ws.Range("H65").Value = 2322.3447
$ws.Range("I65").Value = 2288.5454
$ws.Range("J65").Value = 2428.5715
$ws.Range("K65").Value = 11442.727
$ws.Range("L65").Value = 12142.8575
$ws.Range("M65").Value = -8322.726999999999
$ws.Range("N65").Value = -18382.8575

$ws.Range("H141").Value = 66515.42999999999
$ws.Range("J141").Value = 72552
$ws.Range("L141").Value = 72552
$ws.Range("N141").Value = -82912

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1766.6316
$ws.Range("J109").Value = 3399.5715
$ws.Range("L109").Value = 10198.7145
$ws.Range("N109").Value = -12278.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 29633.334
$ws.Range("J51").Value = 29633.334
$ws.Range("L51").Value = 29633.334
$ws.Range("N51").Value = -30651.334

$ws.Range("H57").Value = 15179.8
$ws.Range("J57").Value = 16156.929
$ws.Range("L57").Value = 16156.929
$ws.Range("N57").Value = -17796.929

$ws.Range("H64").Value = 28300
$ws.Range("J64").Value = 28300
$ws.Range("L64").Value = 28300
$ws.Range("N64").Value = -28796

$ws.Range("H67").Value = 28300
$ws.Range("J67").Value = 28300
$ws.Range("L67").Value = 28300
$ws.Range("N67").Value = -30016

$ws.Range("H113").Value = 1210
$ws.Range("I113").Value = 846.6667
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 846.6667
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 1323.3333
$ws.Range("N113").Value = -6640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2411.4285
$ws.Range("I61").Value = 1720
$ws.Range("J61").Value = 3333.3333
$ws.Range("K61").Value = 1720
$ws.Range("L61").Value = 3333.3333
$ws.Range("M61").Value = -1518
$ws.Range("N61").Value = -3737.3333

$ws.Range("H68").Value = 7850.4
$ws.Range("J68").Value = 3107.8667
$ws.Range("L68").Value = 3107.8667
$ws.Range("N68").Value = -4605.8667

$ws.Range("H71").Value = 7850.4
$ws.Range("J71").Value = 3107.8667
$ws.Range("L71").Value = 15539.3335
$ws.Range("N71").Value = -23027.3335

$ws.Range("H113").Value = 2411.4285
$ws.Range("I113").Value = 1720
$ws.Range("J113").Value = 3333.3333
$ws.Range("K113").Value = 1720
$ws.Range("L113").Value = 3333.3333
$ws.Range("M113").Value = 450
$ws.Range("N113").Value = -7673.3333

$ws.Range("H132").Value = 3896.2856
$ws.Range("I132").Value = 5605.9443
$ws.Range("J132").Value = 2086.0588
$ws.Range("K132").Value = 16817.8329
$ws.Range("L132").Value = 6258.176399999999
$ws.Range("M132").Value = -14287.8329
$ws.Range("N132").Value = -11318.1764
